$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (C, D, E) before the existing "kelas" column,
# pushing it from C to F, then add a new "unit" column in G.
$ws.Range("C1:E1").EntireColumn.Insert()

# Header row
$ws.Range("C1").Value = "Jurusan"
$ws.Range("D1").Value = "Inisial Jurusan"
$ws.Range("E1").Value = "Tahun"

# Data rows
$ws.Range("C2").Value = "web programming"
$ws.Range("G1").Value = "unit"
$ws.Range("D2").Value = "wp"
$ws.Range("E2").Value = 2017
$ws.Range("G2").Value = "wb 1"

$ws.Range("C3").Value = "akuntansi"
$ws.Range("D3").Value = "ak"
$ws.Range("E3").Value = 2017
$ws.Range("G3").Value = "wb 2"

$ws.Range("C4").Value = "software engineer"
$ws.Range("D4").Value = "se"
$ws.Range("E4").Value = 2017
$ws.Range("G4").Value = "wb 1"

$ws.Range("C5").Value = "akuntansi"
$ws.Range("D5").Value = "ak"
$ws.Range("E5").Value = 2017
$ws.Range("G5").Value = "wb 2"

$ws.Range("C6").Value = "software engineer"
$ws.Range("D6").Value = "se"
$ws.Range("E6").Value = 2017
$ws.Range("G6").Value = "wb 1"

$ws.Range("H12").Select()
